$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "disease" entry in column A at row 25 (row 24 intentionally left
# blank, matching the source data gap). Reuse the formatting already applied
# to the other column-A entries (e.g. A23) instead of restyling manually so
# no extra font/style records are introduced.
$ws.Range("A25").Value = "disease"
$ws.Range("A23").Copy() | Out-Null
$ws.Range("A25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update selection to mirror the post-edit cursor position.
$ws.Range("A26").Select() | Out-Null
